$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the famhist_* columns J:Q (merge-import reorder) while keeping
# each header's data aligned with it (values follow their header label).
$ws.Range("J1").Value = "famhist_none"
$ws.Range("K1").Value = "famhist_deaf"
$ws.Range("L1").Value = "famhhist_cardiomyopathy"
$ws.Range("M1").Value = "famhist_encephalopathy"
$ws.Range("N1").Value = "famhist_diabmell"
$ws.Range("O1").Value = "famhist_cardiovasc"
$ws.Range("P1").Value = "famhist_malignancy"
$ws.Range("Q1").Value = "famhist_unknown"

$ws.Range("J4").Value = "Yes"
$ws.Range("K4").Value = "No"
$ws.Range("L4").Value = "No"
$ws.Range("M4").Value = "No"
$ws.Range("N4").Value = "No"
$ws.Range("O4").Value = "No"
$ws.Range("P4").Value = "No"
$ws.Range("Q4").Value = "No"

$ws.Range("J5").Value = "No"
$ws.Range("K5").Value = "No"
$ws.Range("L5").Value = "No"
$ws.Range("M5").Value = "No"
$ws.Range("N5").Value = "No"
$ws.Range("O5").Value = "Yes"
$ws.Range("P5").Value = "Yes"
$ws.Range("Q5").Value = "No"

$ws.Range("J6").Value = "No"
$ws.Range("K6").Value = "No"
$ws.Range("L6").Value = "No"
$ws.Range("M6").Value = "No"
$ws.Range("N6").Value = "No"
$ws.Range("O6").Value = "No"
$ws.Range("P6").Value = "No"
$ws.Range("Q6").Value = "Yes"

# Move the selection to reflect the newly written header range.
$ws.Range("J1:Q1").Select()
